# "wrapping up test file audit"
#
# The optimization_parameters sheet had a stray leftover row (A16:C16 =
# "Sheet", 3, 4) that doesn't belong with the rest of the parameter table
# (row 15 is the Strain/wt/dcin5 header, row 17 is simulation_timepoints).
# Remove that row entirely; everything below shifts up to fill the gap,
# and the now-unused "Sheet" shared string is dropped automatically.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")
$ws.Rows.Item(16).Delete()

# The last-edited/active sheet moves to optimization_diagnostics.
$wb.Worksheets.Item("optimization_diagnostics").Activate()
